# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns for the 7b05ca50-... row (row 5) on both the zh-cn and de-de
# report sheets, and widens the "Error Detail" column to fit the new message.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67dbbb5cd383e3388a0af1ffe49279a5eab71068/e2e/7b05ca50-760b-4c5b-aede-dd1fd27db831.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aafab537c2e706fe7f3b9148e50185fd10766a1/e2e/7b05ca50-760b-4c5b-aede-dd1fd27db831.md."

# Hyperlink style colour (matches the workbook's existing "HyperLink" style: single
# underline, font colour FF6495ED).
$hyperlinkColor = 15570276

function Set-HandbackRow {
    param([string]$SheetName, [string]$XlfFileName, [string]$HandbackDateTime, [string]$HyperlinkTarget)

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File (I5) -> the handed-back markdown file, rendered as a hyperlink
    # just like the other rows in this column.
    $ws.Range("I5").Value = "7b05ca50-760b-4c5b-aede-dd1fd27db831.md"
    $ws.Range("I5").Font.Underline = 2
    $ws.Range("I5").Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($ws.Range("I5"), $HyperlinkTarget, "", "", "7b05ca50-760b-4c5b-aede-dd1fd27db831.md") | Out-Null

    # Latest Handback File (J5)
    $ws.Range("J5").Value = $XlfFileName

    # Latest Handback DateTime (K5)
    $ws.Range("K5").Value = $HandbackDateTime

    # Error Detail (P5)
    $ws.Range("P5").Value = $errorMessage

    # Widen the Error Detail column (P) so the long message is readable.
    $ws.Range("P1").ColumnWidth = 39.17
}

Set-HandbackRow "zh-cn" `
    "7b05ca50-760b-4c5b-aede-dd1fd27db831.afe1eb4272944a2fa5478ac36d666f24a3b085c2.zh-cn.xlf" `
    "2016-10-18 04:22:02" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5e49ae4f03aa1b56e23f60e31f822532543a5878/e2e/7b05ca50-760b-4c5b-aede-dd1fd27db831.md"

Set-HandbackRow "de-de" `
    "7b05ca50-760b-4c5b-aede-dd1fd27db831.afe1eb4272944a2fa5478ac36d666f24a3b085c2.de-de.xlf" `
    "2016-10-18 04:22:31" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5e13839cc52cf5458aed012530029eed6c5d160c/e2e/7b05ca50-760b-4c5b-aede-dd1fd27db831.md"
